# Rotate the "Recorded By" (column G) comma-separated list left by one
# position for every data row that has more than one entry.
# e.g. "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($null -ne $val -and $val -ne "") {
        $parts = $val -split ", "
        if ($parts.Length -gt 1) {
            $rotated = ($parts[1..($parts.Length - 1)] + $parts[0]) -join ", "
            $cell.Value2 = $rotated
        }
    }
}
